# Rename 'Codelists' to 'Cells' and make it the active sheet
# (Close #256)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Codelists")
$ws.Name = "Cells"

# The renamed sheet becomes the active / selected sheet: tabSelected moves
# from "Variables" to "Cells", and the workbook's activeTab index advances.
$ws.Activate()
